# Gantt Chart update: mark the status of the first three sprints in a new
# column F ("Achieved" / "Not completed in this sprint").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = "Achieved"
$ws.Range("F7").Value = "Achieved"
$ws.Range("F8").Value = "Not completed in this sprint"

# Leave the selection where the author left off after typing the new values.
$ws.Range("F9").Select()
